$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Row 2 ---
$ws1.Range("A2").Value = '메뉴 주문 시, 재료 재고량을 반영하기'
$ws1.Range("B2").Value = '테이블에서 메뉴를 주문할 시, 해당 메뉴의 재료가 줄어든다. 재료가 부족하다면 주문할 수 없다.'
$ws1.Range("C2").Value = 43594
$ws1.Range("D2").Value = 43597
$ws1.Range("E2").Value = 'Table 파일 commit'
$ws1.Range("F2").Value = '보안성이 낮음. 개선이 필요함'

# --- Row 3 ---
$ws1.Range("A3").Value = 'Menu 파일의 재료 입력란 수정'
$ws1.Range("B3").Value = '기존에는 한 메뉴에 하나의 재료만 입력할 수 있었음. 수정 뒤엔 여러 메뉴 입력 가능'
$ws1.Range("C3").Value = 43594
$ws1.Range("D3").Value = 43597
$ws1.Range("E3").Value = 'Menu 파일 commit'
$ws1.Range("F3").Value = '재료 문자열을 분리하는 기능을 넣지 못함. 수정 필요.'

# --- Row 4 (new content) ---
$ws1.Range("A4").Value = 'table.java 테이블 추가/제거 관련 오류 수정'
$ws1.Range("B4").Value = '기존에는 테이블 추가/제거가 gui에 반영되지 않음. 그리고 테이블 메뉴 추가/ 결제 기능이 반영되지 않음. 이벤트 핸들러 수정 후, 테이블 추가/제거 기능 사용 가능. 테이블 메뉴 추가/결제 기능 사용 가능'
$ws1.Range("C4").Value = 43600
$ws1.Range("D4").Value = 43606
$ws1.Range("E4").Value = 'Table 파일 commit'
$ws1.Range("F4").Value = '메뉴 선택하는 JComboBox 동기화 방법이 효율적이지 않음. 다른 방법 필요.'

# --- Row 5 (new content) ---
$ws1.Range("A5").Value = 'table.java 기타 오류 수정'
$ws1.Range("B5").Value = '1) 기존에 테이블 미선택 오류, 테이블 삭제 오류, 재료 소진 오류 등이 사용자에게 보여지지 않던 것을 수정함. 2) 메뉴 추가 코드를 개선 '
$ws1.Range("C5").Value = 43600
$ws1.Range("D5").Value = 43606
$ws1.Range("E5").Value = 'Table 파일 commit'
# F5 intentionally left blank

# --- Row heights reflecting the extra content (rows 4 & 5 grew substantially) ---
$ws1.Rows(4).RowHeight = 86.25
$ws1.Rows(5).RowHeight = 61.5

# --- Switch the active sheet/tab from "이정원" (sheet4) to "박서린" (sheet1) ---
$ws1.Activate()
$ws1.Range("B11").Select()
